# Alumni sheet update: replace the 5th alumni row's name with a new
# entry ("Jaya Sagar"). This introduces a brand-new shared string and
# points A5 at it (previously A5 reused the "Aryan Khandal" string).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "Jaya Sagar"

# Restore the view/selection state recorded for the sheet (the cursor
# had moved to F14 with the viewport scrolled so column B is leftmost).
$ws.Range("F14").Select() | Out-Null
